$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.721.90"
$ws.Range("E2").Value = "  +3.03%  "

# Row 3
$ws.Range("D3").Value = "1.864.90"
$ws.Range("E3").Value = "  +3.00%  "

# Row 4
$ws.Range("D4").Value = "'1.037"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.97%  "

# Row 5
$ws.Range("D5").Value = "'324.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.23%  "

# Row 6
$ws.Range("D6").Value = "'1.033"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.61%  "

# Row 7
$ws.Range("D7").Value = "'0.4418"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.01%  "

# Row 8
$ws.Range("D8").Value = "'0.3804"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.03%  "

# Row 9
$ws.Range("D9").Value = "'0.07459"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.97%  "

# Row 10
$ws.Range("D10").Value = "'0.8842"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.46%  "

# Row 11
$ws.Range("D11").Value = "'21.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.88%  "

# Row 12
$ws.Range("D12").Value = "1.880.61"
$ws.Range("E12").Value = "  -8.55%  "

# Row 13
$ws.Range("D13").Value = "'5.562"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.19%  "

# Row 14
$ws.Range("D14").Value = "'6.751"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.85%  "

# Row 15
$ws.Range("D15").Value = "'0.07213"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.19%  "

# Row 16
$ws.Range("D16").Value = "'83.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.83%  "

# Row 17
$ws.Range("D17").Value = "'1.037"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.47%  "

# Row 18
$ws.Range("D18").Value = "'0.000009098"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.07%  "

# Row 19
$ws.Range("D19").Value = "'1.033"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.64%  "

# Row 20
$ws.Range("D20").Value = "'15.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("D21").Value = "27.738.44"
$ws.Range("E21").Value = "  +2.88%  "

# Row 22
$ws.Range("D22").Value = "'5.312"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.22%  "

# Row 23
$ws.Range("D23").Value = "'11.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.35%  "

# Row 24
$ws.Range("D24").Value = "'158.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.88%  "

# Row 25
$ws.Range("D25").Value = "'1.934"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.66%  "

# Row 26
$ws.Range("D26").Value = "'18.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.98%  "

# Row 27
$ws.Range("E27").Value = "  +4.44%  "

# Row 28
$ws.Range("D28").Value = "'5.344"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.09%  "

# Row 29
$ws.Range("D29").Value = "'117.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.50%  "

# Row 30
$ws.Range("D30").Value = "'0.09098"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.77%  "

# Row 31
$ws.Range("E31").Value = "  +5.15%  "

# Row 32
$ws.Range("D32").Value = "'0.7680"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.71%  "

# Row 33
$ws.Range("D33").Value = "'4.575"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.50%  "

# Row 34
$ws.Range("D34").Value = "'2.915"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.90%  "

# Row 35
$ws.Range("D35").Value = "'1.034"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.65%  "

# Row 36
$ws.Range("D36").Value = "'1.160"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.15%  "

# Row 37
$ws.Range("D37").Value = "'0.01990"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.57%  "

# Row 38
$ws.Range("D38").Value = "'0.05348"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.29%  "

# Row 39
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.69%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5205"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.39%  "

# Row 41
$ws.Range("D41").Value = "'0.1694"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.67%  "

# Row 42
$ws.Range("D42").Value = "'6.849"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.13%  "

# Row 43
$ws.Range("D43").Value = "'8.762"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.88%  "

# Row 44
$ws.Range("D44").Value = "'109.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.36%  "

# Row 45
$ws.Range("E45").Value = "  +2.61%  "

# Row 46
$ws.Range("D46").Value = "'1.735"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.19%  "

# Row 47
$ws.Range("D47").Value = "'0.4693"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.06%  "

# Row 48
$ws.Range("E48").Value = "  +2.39%  "

# Row 49
$ws.Range("D49").Value = "'1.868"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "

# Row 50
$ws.Range("D50").Value = "'39.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.80%  "

# Row 51
$ws.Range("D51").Value = "'0.9361"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.81%  "
